# Cucumber default build profile
# Append two new scenario-status rows (56 and 57) to the CampusTest sheet,
# mirroring the existing "Login with valid username and password" /
# "Create Country" PASSED/chrome rows with new timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(56, 1).Value = "Login with valid username and password"
$ws.Cells.Item(56, 2).Value = "PASSED"
$ws.Cells.Item(56, 3).Value = "chrome"
$ws.Cells.Item(56, 4).Value = "11_06_23_154502"

$ws.Cells.Item(57, 1).Value = "Create Country"
$ws.Cells.Item(57, 2).Value = "PASSED"
$ws.Cells.Item(57, 3).Value = "chrome"
$ws.Cells.Item(57, 4).Value = "11_06_23_154510"
